$wb = $excel.ActiveWorkbook

# The "8PFormSheet" (internal codeName "Rooms") is the sheet whose log data
# gets reset for a new day, and it becomes the active tab.
$rooms = $wb.Worksheets.Item("8PFormSheet")

# Clear out the previous day's log entries (rows 2-29, columns A-F),
# keeping row/column formatting intact.
$rooms.Range("A2:F29").ClearContents()

# A couple of cells happen to carry an explicit style that is identical to
# their column's default style; once their value is cleared Excel drops the
# now-redundant style override. Re-apply the default "Normal" style to match.
$rooms.Range("B28").Style = "Normal"
$rooms.Range("A29").Style = "Normal"
$rooms.Range("B29").Style = "Normal"

# Log the new entry for the next day.
$rooms.Range("A2").Value = 43068
$rooms.Range("B2").Value = "MS"
$rooms.Range("C2").Value = "NO"
$rooms.Range("D2").Value = "2"

# Make the 8PFormSheet the active/selected tab.
$rooms.Activate()
